$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RiskList")
$tbl = $ws.ListObjects.Item("RiskList")

# ---------------------------------------------------------------------------
# 1) Add a brand new risk row (R6) to the RiskList table
# ---------------------------------------------------------------------------
$newRow = $tbl.ListRows.Add()

# Copy the formatting of the previous data row (R5 / sheet row 8) onto the
# freshly inserted row so it keeps matching borders/fonts/number formats.
$ws.Range("B8:G8").Copy()
$ws.Range("B9:G9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B9").Formula = '="R"&ROW($A6)'
$ws.Range("C9").Value2 = "Způsob načítání dat"
$ws.Range("D9").Value2 = "Špatný způsob, jak data načítat povede k složité udržitelnosti a následovném rozšiřování databáze."

# ---------------------------------------------------------------------------
# 2) Row 5 (risk R2): "Neznalost architektury" -> "Neznalost technologie"
#    (text of the risk / impact / mitigation action updated)
# ---------------------------------------------------------------------------
$ws.Range("C5").Value2 = "Neznalost technologie"
$ws.Range("D5").Value2 = "Bez známosti použité technologie není možné začít vývoj projektu."
$ws.Range("E5").Value2 = "Porovnat vhodné technologie, vybrat tu, která splňuje nejvíc požadavků. Seznámit se s jejími součásmi vytvořením drobného prototypu."

# ---------------------------------------------------------------------------
# 3) Finish filling in the new risk row (R6)
# ---------------------------------------------------------------------------
$ws.Range("E9").Value2 = "Prozkoumat nejpoužívanější způsoby správy dat pro vybranou technologii. Zvolit tu pro projekt nejvhodnější."
$ws.Range("F9").Value2 = 3
$ws.Range("G9").Value2 = 0.75

# ---------------------------------------------------------------------------
# 4) Row 6 (risk R3): Priorita (priority) changed from 2 to 3
# ---------------------------------------------------------------------------
$ws.Range("F6").Value2 = 3

# ---------------------------------------------------------------------------
# 5) Update view/selection state to match the author's last action
# ---------------------------------------------------------------------------
$ws.Range("H9").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
